$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin -> Bitcoin
$ws.Range("D2").Value = "30.647.16"
$ws.Range("E2").Value = "  -1.20%  "

# Row 3: Ethereum -> Ethereum
$ws.Range("D3").Value = "1.919.25"
$ws.Range("E3").Value = "  -1.95%  "

# Row 4: TetherUSD -> TetherUSD
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.51%  "

# Row 5: BNB -> BNB
$ws.Range("D5").Value = "'238.17"
$ws.Range("E5").Value = "  -2.79%  "

# Row 6: USDC -> USDC
$ws.Range("D6").Value = "'1.005"
$ws.Range("E6").Value = "  +0.39%  "

# Row 7: XRP -> XRP
$ws.Range("D7").Value = "'0.4770"
$ws.Range("E7").Value = "  -1.95%  "

# Row 8: Cardano -> Cardano
$ws.Range("D8").Value = "'0.2873"
$ws.Range("E8").Value = "  -2.75%  "

# Row 9: Dogecoin -> Dogecoin
$ws.Range("D9").Value = "'0.06668"
$ws.Range("E9").Value = "  -2.15%  "

# Row 10: Solana -> Solana
$ws.Range("D10").Value = "'18.69"
$ws.Range("E10").Value = "  -2.45%  "

# Row 11: Litecoin -> Litecoin
$ws.Range("D11").Value = "'103.51"
$ws.Range("E11").Value = "  -3.44%  "

# Row 12: TRON -> WrappedEther
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.929.52"
$ws.Range("E12").Value = "  -1.51%  "

# Row 13: WrappedEther -> TRON
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.07716"
$ws.Range("E13").Value = "  -1.11%  "

# Row 14: Polkadot -> Polkadot
$ws.Range("D14").Value = "'5.221"
$ws.Range("E14").Value = "  -4.20%  "

# Row 15: Polygon -> Polygon
$ws.Range("D15").Value = "'0.6841"
$ws.Range("E15").Value = "  -2.63%  "

# Row 16: BitcoinCash -> BitcoinCash
$ws.Range("D16").Value = "'264.89"
$ws.Range("E16").Value = "  -6.53%  "

# Row 17: WrappedBTC -> WrappedBTC
$ws.Range("D17").Value = "30.679.78"
$ws.Range("E17").Value = "  -1.11%  "

# Row 18: ShibaInu -> Dai
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "'1.002"
$ws.Range("E18").Value = "  +0.15%  "

# Row 19: Dai -> ShibaInu
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000007492"
$ws.Range("E19").Value = "  -2.43%  "

# Row 20: Avalanche -> Avalanche
$ws.Range("D20").Value = "'12.69"
$ws.Range("E20").Value = "  -3.71%  "

# Row 21: Uniswap -> Uniswap
$ws.Range("D21").Value = "'5.449"
$ws.Range("E21").Value = "  -0.94%  "

# Row 22: BinanceUSD -> BinanceUSD
$ws.Range("D22").Value = "'1.008"
$ws.Range("E22").Value = "  +0.71%  "

# Row 23: BitDAO -> BitDAO
$ws.Range("D23").Value = "'0.4555"
$ws.Range("E23").Value = "  -8.95%  "

# Row 24: Chainlink -> Chainlink
$ws.Range("D24").Value = "'6.337"
$ws.Range("E24").Value = "  -2.32%  "

# Row 25: Cosmos -> Cosmos
$ws.Range("D25").Value = "'9.676"
$ws.Range("E25").Value = "  -1.16%  "

# Row 26: Monero -> Monero
$ws.Range("D26").Value = "'163.13"
$ws.Range("E26").Value = "  -4.49%  "

# Row 27: EthereumClassic -> EthereumClassic
$ws.Range("D27").Value = "'18.90"
$ws.Range("E27").Value = "  -5.38%  "

# Row 28: LidoDAOToken -> LidoDAOToken
$ws.Range("D28").Value = "'2.098"
$ws.Range("E28").Value = "  -5.18%  "

# Row 29: Stellar -> Toncoin
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'1.397"
$ws.Range("E29").Value = "  -0.91%  "

# Row 30: Toncoin -> Stellar
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "'0.1017"
$ws.Range("E30").Value = "  -3.86%  "

# Row 31: Filecoin -> PancakeSwap
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.528"
$ws.Range("E31").Value = "  -3.35%  "

# Row 32: PancakeSwap -> Filecoin
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.392"
$ws.Range("E32").Value = "  -4.48%  "

# Row 33: InternetComputer(DFINITY) -> InternetComputer(DFINITY)
$ws.Range("D33").Value = "'4.223"
$ws.Range("E33").Value = "  -5.15%  "

# Row 34: Hedera -> Hedera
$ws.Range("D34").Value = "'0.04748"
$ws.Range("E34").Value = "  -3.62%  "

# Row 35: ImmutableX -> ImmutableX
$ws.Range("D35").Value = "'0.7324"
$ws.Range("E35").Value = "  -3.64%  "

# Row 36: ARBITRUM -> ARBITRUM
$ws.Range("D36").Value = "'1.120"
$ws.Range("E36").Value = "  -4.24%  "

# Row 37: Frax -> Frax
$ws.Range("D37").Value = "'1.002"
$ws.Range("E37").Value = "  +0.18%  "

# Row 38: HuobiToken -> HuobiToken
$ws.Range("D38").Value = "'2.742"
$ws.Range("E38").Value = "  +0.47%  "

# Row 39: VeChain -> VeChain
$ws.Range("D39").Value = "'0.01949"
$ws.Range("E39").Value = "  -3.15%  "

# Row 40: MXToken -> MXToken
$ws.Range("D40").Value = "'2.654"
$ws.Range("E40").Value = "  -1.78%  "

# Row 41: FraxShare -> FraxShare
$ws.Range("D41").Value = "'6.317"
$ws.Range("E41").Value = "  -3.09%  "

# Row 42: Aave -> Aave
$ws.Range("D42").Value = "'75.44"
$ws.Range("E42").Value = "  -1.44%  "

# Row 43: RenderToken -> RenderToken
$ws.Range("D43").Value = "'2.010"
$ws.Range("E43").Value = "  -4.88%  "

# Row 44: TrustWalletToken -> TrustWalletToken
$ws.Range("D44").Value = "'0.8673"
$ws.Range("E44").Value = "  -2.26%  "

# Row 45: TheSandbox -> Quant
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'106.36"
$ws.Range("E45").Value = "  -2.61%  "

# Row 46: Quant -> TheSandbox
$ws.Range("B46").Value = "TheSandbox"
$ws.Range("C46").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D46").Value = "'0.4285"
$ws.Range("E46").Value = "  -3.95%  "

# Row 47: PaxDollar -> PaxDollar
$ws.Range("D47").Value = "'1.004"
$ws.Range("E47").Value = "  +0.30%  "

# Row 48: Aptos -> Aptos
$ws.Range("D48").Value = "'7.567"
$ws.Range("E48").Value = "  -7.30%  "

# Row 49: Maker -> Maker
$ws.Range("D49").Value = "'959.10"
$ws.Range("E49").Value = "  -3.78%  "

# Row 50: Algorand -> Algorand
$ws.Range("D50").Value = "'0.1199"
$ws.Range("E50").Value = "  -4.75%  "

# Row 51: Elrond -> Elrond
$ws.Range("D51").Value = "'35.10"
$ws.Range("E51").Value = "  -1.93%  "
